$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.351.98"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -2.61%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.856.46"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -2.64%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.09%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'327.60"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.11%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.09%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4549"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -2.37%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3906"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -2.17%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'47.82"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -10.13%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07918"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -5.37%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'1.012"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.65%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -2.51%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.857.08"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.46%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.924"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -1.86%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'7.163"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -3.24%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'1.002"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.09%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = "'0.06650"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +1.10%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = "'86.15"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -3.42%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.00001029"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.03%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'17.19"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -3.81%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'1.002"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.01%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.505"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -3.35%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'27.351.64"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -2.73%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'10.89"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -3.44%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.288"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.57%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.082.81"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -3.83%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'154.25"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.32%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'19.99"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.41%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.065"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -2.22%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'5.466"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.77%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'121.23"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.89%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.9510"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -1.63%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.09372"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.66%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.447"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +1.55%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'3.590"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.32%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'5.261"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -4.39%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.06052"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.07%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.02228"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -2.71%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'1.217"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.06%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'8.099"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -7.90%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.02%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.5929"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -2.80%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.1889"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.00%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'10.14"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -7.19%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.59%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.5614"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -3.37%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'12.09"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -4.40%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'3.389"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -1.42%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.919"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -4.43%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.06748"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.29%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'108.46"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.54%  "
$ws.Range('E51').Style = 'Normal'
